# nová verze 2.0, možnosti mountění, oprava chyb

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------
# Sheet "ip_address_list" (sheet1) - fix/cleanup a few rows
# -------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")

# row1: drop the stray "None" note
$ws1.Range("D1").ClearContents()

# row2: id was garbage ("5145566") -> "dsff", note was "None" -> long "f..." string
$ws1.Range("A2").Value = "dsff"
$ws1.Range("D2").Value = "fffffffffffffffffffffffffffffff"

# row3: id "dsff" -> "5" (kept as text), drop the long note that moved to row2
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "5"
$ws1.Range("D3").ClearContents()

# row4: id "5" -> "514" (kept as text), fix the IP typo 100 -> 14
$ws1.Range("A4").NumberFormat = "@"
$ws1.Range("A4").Value = "514"
$ws1.Range("B4").Value = "192.168.14.241"

# -------------------------------------------------------------
# Sheet "disc_list" (sheet2) - new layout with mount info + hyperlink
# -------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("disc_list")

# old row1 is being fully replaced by new content; clear any leftover
# per-cell formatting (e.g. the old B1 had a custom font) before writing.
$ws2.Range("A1:F1").ClearFormats()

$ws2.Range("A1").NumberFormat = "@"
$ws2.Range("A1").Value = "518"
$ws2.Range("B1").Value = "V"
$ws2.Range("C1").Value = "\\192.168.208.200\10_vision"
$ws2.Range("D1").Value = "jhv_vision"
$ws2.Range("E1").Value = "Jhv*2708"
$ws2.Range("F1").Value = "první sít, ixon"

$ws2.Range("A2").Value = "514-2"
$ws2.Range("B2").Value = "T"
$ws2.Range("C2").Value = "\\192.168.14.245\Data\Kamery"
$ws2.Range("D2").Value = "Vision"
$ws2.Range("E2").Value = "*Jhv2708"

$ws2.Hyperlinks.Add($ws2.Range("C1"), "\\192.168.208.200\10_vision") | Out-Null

$ws2.Columns.Item(2).ColumnWidth = 8.42578125
$ws2.Columns.Item(3).ColumnWidth = 32

$ws2.Range("C1").Select()

Write-Host "edits applied"
